# The underlying edit (per the commit diff) is a re-shuffle of the row
# content for the species-observation records on the "Artfynd" sheet: row N
# ends up holding what used to be in row Map[N] (rows 18-19 swap; rows
# 21-29 cycle). Row 20 is untouched. Only the columns that actually differ
# between the paired rows are touched (Id/A, Taxonsorteringsordning/B,
# Rodlistade/D, TaxonId/E, Artnamn/F, Vetenskapligt namn/G, Auktor/H,
# Ost/Q, Nord/R) - every other column (dates, empty placeholder cells,
# location, etc.) is identical between the paired rows so there is nothing
# to write there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = "A","B","D","E","F","G","H","Q","R"

$map = @{
    18 = 19
    19 = 18
    21 = 26
    22 = 29
    23 = 24
    24 = 22
    25 = 28
    26 = 25
    27 = 21
    28 = 27
    29 = 23
}

# Snapshot the current (pre-edit) values of the columns that matter, for
# every row involved, before any writes happen - needed because the
# remapping above is cyclic (e.g. 21 -> 26 -> 25 -> 28 -> 27 -> 21), so a
# row's original content must be captured before it is overwritten.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range($col + $r).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value2 = $srcVals[$col]
    }
}
